# Insert a new weekly price record as row 73 in the sheet, pushing the
# existing row 73 (and everything after it) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(73).Insert()

$ws.Cells.Item(73, 1).Value  = 11
$ws.Cells.Item(73, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value  = "Bíobío"
$ws.Cells.Item(73, 4).Value  = 44546
$ws.Cells.Item(73, 5).Value  = 8
$ws.Cells.Item(73, 6).Value  = 100112023
$ws.Cells.Item(73, 7).Value  = "Brócoli"
$ws.Cells.Item(73, 8).Value  = "Sin especificar"
$ws.Cells.Item(73, 9).Value  = "Primera"
$ws.Cells.Item(73, 10).Value = 2700
$ws.Cells.Item(73, 11).Value = 500
$ws.Cells.Item(73, 12).Value = 600
$ws.Cells.Item(73, 13).Value = 544
$ws.Cells.Item(73, 14).Value = "`$/unidad"
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 544
$ws.Cells.Item(73, 17).Value = 1
$ws.Cells.Item(73, 18).Value = "Hortaliza"
